# Scheduled-runner refresh of market price snapshots (currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ and the derived Leve profit
# columns H-N) across the eight class sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1790.72
$ws.Range("I100").Value = 1271.875
$ws.Range("K100").Value = 1271.875
$ws.Range("M100").Value = -730.875
$ws.Range("H103").Value = 1499.9445
$ws.Range("I103").Value = 1407.6154
$ws.Range("J103").Value = 1740
$ws.Range("K103").Value = 4222.8462
$ws.Range("L103").Value = 5220
$ws.Range("M103").Value = -3636.8462
$ws.Range("N103").Value = -6392
$ws.Range("H137").Value = 3057.795
$ws.Range("I137").Value = 2071.0344
$ws.Range("J137").Value = 5919.4
$ws.Range("K137").Value = 6213.1032
$ws.Range("L137").Value = 17758.2
$ws.Range("M137").Value = -3663.1032
$ws.Range("N137").Value = -22858.2
$ws.Range("H138").Value = 2379.0632
$ws.Range("I138").Value = 1251.8
$ws.Range("J138").Value = 3069.2246
$ws.Range("K138").Value = 3755.4
$ws.Range("L138").Value = 9207.6738
$ws.Range("M138").Value = 1384.6
$ws.Range("N138").Value = -19487.6738

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8932254
$ws.Range("I32").Value = 9618428
$ws.Range("K32").Value = 9618428
$ws.Range("M32").Value = -9618141
$ws.Range("H80").Value = 79912
$ws.Range("J80").Value = 79912
$ws.Range("L80").Value = 79912
$ws.Range("N80").Value = -81908
$ws.Range("H83").Value = 79912
$ws.Range("J83").Value = 79912
$ws.Range("L83").Value = 239736
$ws.Range("N83").Value = -249720
$ws.Range("H93").Value = 30000
$ws.Range("J93").Value = 30000
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -34992
$ws.Range("H132").Value = 3218.3333
$ws.Range("I132").Value = 1476.9565
$ws.Range("K132").Value = 4430.8695
$ws.Range("M132").Value = -1900.8695

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H93").Value = 69996.336
$ws.Range("J93").Value = 69996.336
$ws.Range("L93").Value = 69996.336
$ws.Range("N93").Value = -73740.336
$ws.Range("H97").Value = 20872.182
$ws.Range("I97").Value = 10578.125
$ws.Range("J97").Value = 48323
$ws.Range("K97").Value = 10578.125
$ws.Range("L97").Value = 48323
$ws.Range("M97").Value = -9587.125
$ws.Range("N97").Value = -50305

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 718.8333
$ws.Range("I16").Value = 602.36365
$ws.Range("K16").Value = 602.36365
$ws.Range("M16").Value = -315.36365
$ws.Range("H76").Value = 8886
$ws.Range("I76").Value = 8886
$ws.Range("K76").Value = 8886
$ws.Range("M76").Value = -8571
$ws.Range("H79").Value = 8886
$ws.Range("I79").Value = 8886
$ws.Range("K79").Value = 8886
$ws.Range("M79").Value = -7794
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").ClearContents()
$ws.Range("N97").Value = 0
$ws.Range("H99").Value = 5988.6
$ws.Range("I99").Value = 6089.125
$ws.Range("J99").Value = 5586.5
$ws.Range("K99").Value = 6089.125
$ws.Range("L99").Value = 5586.5
$ws.Range("M99").Value = -4591.125
$ws.Range("N99").Value = -8582.5
$ws.Range("H113").Value = 718.8333
$ws.Range("I113").Value = 602.36365
$ws.Range("K113").Value = 602.36365
$ws.Range("M113").Value = 1567.63635
$ws.Range("H118").Value = 87989
$ws.Range("J118").Value = 87989
$ws.Range("L118").Value = 87989
$ws.Range("N118").Value = -91303
$ws.Range("H122").Value = 862.2941
$ws.Range("I122").Value = 903.0833
$ws.Range("J122").Value = 764.4
$ws.Range("K122").Value = 2709.2499
$ws.Range("L122").Value = 2293.2
$ws.Range("M122").Value = -259.2498999999998
$ws.Range("N122").Value = -7193.2
$ws.Range("H126").Value = 5988.6
$ws.Range("I126").Value = 6089.125
$ws.Range("J126").Value = 5586.5
$ws.Range("K126").Value = 18267.375
$ws.Range("L126").Value = 16759.5
$ws.Range("M126").Value = -15797.375
$ws.Range("N126").Value = -21699.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 9481.333000000001
$ws.Range("I82").Value = 8444
$ws.Range("J82").Value = 10000
$ws.Range("K82").Value = 25332
$ws.Range("L82").Value = 30000
$ws.Range("M82").Value = -24926
$ws.Range("N82").Value = -30812
$ws.Range("H85").Value = 9481.333000000001
$ws.Range("I85").Value = 8444
$ws.Range("J85").Value = 10000
$ws.Range("K85").Value = 25332
$ws.Range("L85").Value = 30000
$ws.Range("M85").Value = -23928
$ws.Range("N85").Value = -32808
$ws.Range("H125").Value = 3636.182

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 25000
$ws.Range("I33").Value = 20000
$ws.Range("J33").Value = 30000
$ws.Range("K33").Value = 20000
$ws.Range("L33").Value = 30000
$ws.Range("M33").Value = -19748
$ws.Range("N33").Value = -30504
$ws.Range("H38").Value = 5000
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H40").Value = 24578
$ws.Range("J40").Value = 24578
$ws.Range("L40").Value = 24578
$ws.Range("N40").Value = -24880
$ws.Range("H80").Value = 3750.375
$ws.Range("I80").Value = 3600.8
$ws.Range("K80").Value = 3600.8
$ws.Range("M80").Value = -2602.8
$ws.Range("H83").Value = 3750.375
$ws.Range("I83").Value = 3600.8
$ws.Range("K83").Value = 18004
$ws.Range("M83").Value = -13012
$ws.Range("H92").Value = 19724.7
$ws.Range("J92").Value = 19724.7
$ws.Range("L92").Value = 19724.7
$ws.Range("N92").Value = -23468.7
$ws.Range("H93").Value = 57999.5
$ws.Range("J93").Value = 57999.5
$ws.Range("L93").Value = 57999.5
$ws.Range("N93").Value = -61743.5
$ws.Range("H97").Value = 1084.1818
$ws.Range("I97").Value = 1171.55
$ws.Range("J97").Value = 210.5
$ws.Range("K97").Value = 1171.55
$ws.Range("L97").Value = 210.5
$ws.Range("M97").Value = -675.55
$ws.Range("N97").Value = -1202.5
$ws.Range("H102").Value = 2887.9487
$ws.Range("I102").Value = 2352.162
$ws.Range("J102").Value = 12800
$ws.Range("K102").Value = 2352.162
$ws.Range("L102").Value = 12800
$ws.Range("M102").Value = -730.1619999999998
$ws.Range("N102").Value = -16044
$ws.Range("H132").Value = 50007630
$ws.Range("I132").Value = 71430184
$ws.Range("J132").Value = 21683.834
$ws.Range("K132").Value = 214290552
$ws.Range("L132").Value = 65051.50199999999
$ws.Range("M132").Value = -214288022
$ws.Range("N132").Value = -70111.50199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 116656
$ws.Range("I7").Value = 5998
$ws.Range("J7").Value = 130488.25
$ws.Range("K7").Value = 5998
$ws.Range("L7").Value = 130488.25
$ws.Range("M7").Value = -5886
$ws.Range("N7").Value = -130712.25
$ws.Range("H40").Value = 3567.25
$ws.Range("I40").Value = 2756.3333
$ws.Range("J40").Value = 6000
$ws.Range("K40").Value = 2756.3333
$ws.Range("L40").Value = 6000
$ws.Range("M40").Value = -2620.3333
$ws.Range("N40").Value = -6272
$ws.Range("H46").Value = 5795.2144
$ws.Range("J46").Value = 10716.667
$ws.Range("L46").Value = 10716.667
$ws.Range("N46").Value = -11092.667
$ws.Range("H55").Value = 55556020
$ws.Range("I55").Value = 62500430
$ws.Range("J55").Value = 743.5
$ws.Range("K55").Value = 62500430
$ws.Range("L55").Value = 743.5
$ws.Range("M55").Value = -62500257
$ws.Range("N55").Value = -1089.5
$ws.Range("H61").Value = 1998.3334
$ws.Range("I61").Value = 1997.5
$ws.Range("K61").Value = 1997.5
$ws.Range("M61").Value = -1795.5
$ws.Range("H113").Value = 1998.3334
$ws.Range("I113").Value = 1997.5
$ws.Range("K113").Value = 1997.5
$ws.Range("M113").Value = 172.5
$ws.Range("H122").Value = 4498.971
$ws.Range("I122").Value = 4495.8125
$ws.Range("K122").Value = 13487.4375
$ws.Range("M122").Value = -11037.4375
$ws.Range("H126").Value = 116656
$ws.Range("I126").Value = 5998
$ws.Range("J126").Value = 130488.25
$ws.Range("K126").Value = 17994
$ws.Range("L126").Value = 391464.75
$ws.Range("M126").Value = -15524
$ws.Range("N126").Value = -396404.75
$ws.Range("H136").Value = 8527.4
$ws.Range("I136").Value = 8975.462
$ws.Range("J136").Value = 8042
$ws.Range("K136").Value = 26926.386
$ws.Range("L136").Value = 24126
$ws.Range("M136").Value = -24376.386
$ws.Range("N136").Value = -29226

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 62028
$ws.Range("J40").Value = 62028
$ws.Range("L40").Value = 62028
$ws.Range("N40").Value = -62326
$ws.Range("H121").Value = 23500
$ws.Range("J121").Value = 23500
$ws.Range("L121").Value = 23500
$ws.Range("N121").Value = -26994
$ws.Range("H132").Value = 2233.3845
$ws.Range("I132").Value = 2355.1738
$ws.Range("K132").Value = 7065.5214
$ws.Range("M132").Value = -4535.5214
